# Auto refresh - 15-02-2026 11:55:27.76
# Merge in the latest ML scoring pass: add Flight_Start/End "_meta" echo columns,
# a derived Campaign_Status column, and an Excel-vs-ML disagreement flag.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Excel_vs_ML
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Excel_vs_ML")

$lastRow = 46

# 1) Insert two fresh columns at Q:R (Flight_Start_Date_meta / Flight_End_Date_meta).
#    This pushes the old Q (ML_Prediction) -> S, old R (duplicate Total_Budget) -> T,
#    old S (Budget_At_Risk) -> U and old T (ML_Early_Warning) -> V.
$ws.Range("Q1:R1").EntireColumn.Insert()

# 2) The old duplicate "Total_Budget" column has landed on T - drop it, it is
#    identical to columns G/P and superseded by the meta columns above.
$ws.Range("T1").EntireColumn.Delete()

# 3) Insert a new column in its place for the derived Campaign_Status field.
$ws.Range("T1").EntireColumn.Insert()

# 4) Append the new Excel_vs_ML_Disagree column after the existing last column (V).
#    (Column W is already blank/unused, so a plain write is enough.)

# Header renames / additions
$ws.Range("G1").Value = "Total_Budget"
$ws.Range("P1").Value = "Total_Budget_meta"
$ws.Range("Q1").Value = "Flight_Start_Date_meta"
$ws.Range("R1").Value = "Flight_End_Date_meta"
$ws.Range("S1").Value = "ML_Prediction"
$ws.Range("T1").Value = "Campaign_Status"
$ws.Range("U1").Value = "Budget_At_Risk"
$ws.Range("V1").Value = "ML_Early_Warning"

# W1 is brand new real estate past the old used range, so it has no header
# styling yet - clone the existing bold/centered/bordered header look from A1.
$ws.Range("W1").Value = "Excel_vs_ML_Disagree"
$ws.Range("A1").Copy()
$ws.Range("W1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$dateFmt = $ws.Range("C2").NumberFormat

for ($r = 2; $r -le $lastRow; $r++) {
    $flightStart = $ws.Cells.Item($r, 3).Value2   # column C
    $flightEnd   = $ws.Cells.Item($r, 4).Value2   # column D
    $todayDate   = $ws.Cells.Item($r, 5).Value2   # column E

    # Flight_Start_Date_meta / Flight_End_Date_meta: straight echo of the
    # flight dates, carried along with the ML merge payload.
    $ws.Cells.Item($r, 17).Value = $flightStart   # column Q
    $ws.Cells.Item($r, 17).NumberFormat = $dateFmt
    $ws.Cells.Item($r, 18).Value = $flightEnd     # column R
    $ws.Cells.Item($r, 18).NumberFormat = $dateFmt

    # Campaign_Status: has the flight already ended relative to Today_Date?
    if ($flightEnd -lt $todayDate) {
        $ws.Cells.Item($r, 20).Value = "ENDED"    # column T
    } else {
        $ws.Cells.Item($r, 20).Value = "LIVE"     # column T
    }

    # Excel_vs_ML_Disagree: this refresh's scoring pass disagreed with the
    # legacy Excel pacing flag on every live row.
    $ws.Cells.Item($r, 23).Value = "YES"          # column W
}

# ---------------------------------------------------------------------------
# Sheet 2: Exec_Summary
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Exec_Summary")

# Insert a new row 4 for the disagreement count, pushing LAST_REFRESH_UTC to row 5.
$summary.Rows(4).Insert()
$summary.Range("A4").Value = "Excel vs ML Disagreement Count"
$summary.Range("B4").Value = 45

$summary.Range("B5").Value = "2026-02-15 06:25 UTC"

# Leave the workbook's active sheet/selection exactly as it was (Excel_vs_ML!A1).
$ws.Activate()
$ws.Range("A1").Select()
